$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "TASK1509262"
$ws.Range("B2").Value = "Rafael Goncalves Reis"
$ws.Range("C2").Value = "DSS - Brazil - Rio de Janeiro / Cidade Nova"

$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4163)
